$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.8198218995343602
$ws.Range("E2").Value = 0.8198218995343602

# Row 3
$ws.Range("D3").Value = 0.0008712521381848836
$ws.Range("E3").Value = 0.0008712521381848836

# Row 4
$ws.Range("D4").Value = 0.000003041250910198796
$ws.Range("E4").Value = 0.000003041250910198796

# Row 5
$ws.Range("D5").Value = 0.08154112952028401
$ws.Range("E5").Value = 0.08154112952028401

# Row 6
$ws.Range("D6").Value = 0.9335990326873733
$ws.Range("E6").Value = 0.9335990326873733

# Row 7 (Success flips from TRUE to FALSE)
$ws.Range("C7").Value = $false
$ws.Range("D7").Value = 0.000161088601295737
$ws.Range("E7").Value = 0.9998389113987043

# Row 8
$ws.Range("D8").Value = 0.9999999999977665
$ws.Range("E8").Value = 0.00000000000223354668094089

# Row 9
$ws.Range("D9").Value = 0.7977890430993961
$ws.Range("E9").Value = 0.2022109569006039

# Row 10
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 0

# Row 11
$ws.Range("D11").Value = 0.9740772405200334
$ws.Range("E11").Value = 0.02592275947996658
$ws.Range("F11").Value = 1.349751710891724
$ws.Range("G11").Value = 0.7

# Row 12
$ws.Range("D12").Value = 0.0001040467123140103
$ws.Range("E12").Value = 0.0001040467123140103

# Row 13
$ws.Range("D13").Value = 0.000258174366474553
$ws.Range("E13").Value = 0.000258174366474553

# Row 14
$ws.Range("D14").Value = 0.0000000000002387051485909284
$ws.Range("E14").Value = 0.0000000000002387051485909284

# Row 15
$ws.Range("D15").Value = 0.03517894395232058
$ws.Range("E15").Value = 0.03517894395232058

# Row 16
$ws.Range("D16").Value = 0.9960015487859081
$ws.Range("E16").Value = 0.9960015487859081

# Row 17 (Success flips from TRUE to FALSE)
$ws.Range("C17").Value = $false
$ws.Range("D17").Value = 0.0000006264109280367846
$ws.Range("E17").Value = 0.999999373589072

# Row 18
$ws.Range("D18").Value = 0.9999999999999967
$ws.Range("E18").Value = 0.00000000000000333066907387547

# Row 19 (Success flips from TRUE to FALSE)
$ws.Range("C19").Value = $false
$ws.Range("D19").Value = 0.00007080365868497946
$ws.Range("E19").Value = 0.9999291963413151

# Row 20
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 0

# Row 21
$ws.Range("D21").Value = 0.9999585299564866
$ws.Range("E21").Value = 0.000041470043513403
$ws.Range("F21").Value = 2.939692258834839
$ws.Range("G21").Value = 0.7
